$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10.75323438732901 ; $ws.Cells.Item(2, 3).Value = 4.771289220509138 ; $ws.Cells.Item(2, 4).Value = 14.96616765246449 ; $ws.Cells.Item(2, 5).Value = 16.38130147518428 ; $ws.Cells.Item(2, 7).Value = 36.0854662352329 ; $ws.Cells.Item(2, 8).Value = 16.19390760500136 ; $ws.Cells.Item(2, 10).Value = 9.350473265601769 ; $ws.Cells.Item(2, 11).Value = 10.07388184593748 ; $ws.Cells.Item(2, 15).Value = 25.61415476017342
$ws.Cells.Item(3, 2).Value = 10.46477751452044 ; $ws.Cells.Item(3, 3).Value = 4.572673318459346 ; $ws.Cells.Item(3, 4).Value = 14.90269282907347 ; $ws.Cells.Item(3, 5).Value = 16.31791882411145 ; $ws.Cells.Item(3, 7).Value = 36.18443169714828 ; $ws.Cells.Item(3, 8).Value = 16.2477179293135 ; $ws.Cells.Item(3, 10).Value = 9.357354326951455 ; $ws.Cells.Item(3, 11).Value = 9.870793124858572 ; $ws.Cells.Item(3, 15).Value = 25.70194422641802
$ws.Cells.Item(4, 2).Value = 10.28522567712254 ; $ws.Cells.Item(4, 3).Value = 4.445484807901144 ; $ws.Cells.Item(4, 4).Value = 14.86695157755296 ; $ws.Cells.Item(4, 5).Value = 16.28257435831664 ; $ws.Cells.Item(4, 7).Value = 36.25579065227989 ; $ws.Cells.Item(4, 8).Value = 16.28336136712342 ; $ws.Cells.Item(4, 10).Value = 9.362971692164985 ; $ws.Cells.Item(4, 11).Value = 9.745280013545157 ; $ws.Cells.Item(4, 15).Value = 25.76119623819798
$ws.Cells.Item(5, 2).Value = 10.21155469539852 ; $ws.Cells.Item(5, 3).Value = 4.392386596967213 ; $ws.Cells.Item(5, 4).Value = 14.85321101189618 ; $ws.Cells.Item(5, 5).Value = 16.26908056395254 ; $ws.Cells.Item(5, 7).Value = 36.28752288956449 ; $ws.Cells.Item(5, 8).Value = 16.29854089919963 ; $ws.Cells.Item(5, 10).Value = 9.365611231871283 ; $ws.Cells.Item(5, 11).Value = 9.693994760640896 ; $ws.Cells.Item(5, 15).Value = 25.7866842458633
$ws.Cells.Item(6, 2).Value = 10.19929471747769 ; $ws.Cells.Item(6, 3).Value = 4.383494677909749 ; $ws.Cells.Item(6, 4).Value = 14.85097949919252 ; $ws.Cells.Item(6, 5).Value = 16.26689515122987 ; $ws.Cells.Item(6, 7).Value = 36.29295189415686 ; $ws.Cells.Item(6, 8).Value = 16.30110097663338 ; $ws.Cells.Item(6, 10).Value = 9.366070696860996 ; $ws.Cells.Item(6, 11).Value = 9.685472579639075 ; $ws.Cells.Item(6, 15).Value = 25.79099750346484
$ws.Cells.Item(7, 2).Value = 10.28423400661879 ; $ws.Cells.Item(7, 3).Value = 4.444773770305635 ; $ws.Cells.Item(7, 4).Value = 14.8667629156745 ; $ws.Cells.Item(7, 5).Value = 16.28238868065919 ; $ws.Cells.Item(7, 7).Value = 36.25620787806669 ; $ws.Cells.Item(7, 8).Value = 16.28356343362785 ; $ws.Cells.Item(7, 10).Value = 9.363005870761045 ; $ws.Cells.Item(7, 11).Value = 9.744588829422481 ; $ws.Cells.Item(7, 15).Value = 25.76153454697317
$ws.Cells.Item(8, 2).Value = 10.65434512843992 ; $ws.Cells.Item(8, 3).Value = 4.703920899564032 ; $ws.Cells.Item(8, 4).Value = 14.94361779936759 ; $ws.Cells.Item(8, 5).Value = 16.35871253285543 ; $ws.Cells.Item(8, 7).Value = 36.11738440745548 ; $ws.Cells.Item(8, 8).Value = 16.21192090672318 ; $ws.Cells.Item(8, 10).Value = 9.352557069739941 ; $ws.Cells.Item(8, 11).Value = 10.00406605852346 ; $ws.Cells.Item(8, 15).Value = 25.6433127974902
$ws.Cells.Item(9, 2).Value = 11.35623787463041 ; $ws.Cells.Item(9, 3).Value = 5.168723923973403 ; $ws.Cells.Item(9, 4).Value = 15.11942811778435 ; $ws.Cells.Item(9, 5).Value = 16.5361989828056 ; $ws.Cells.Item(9, 7).Value = 35.92966639874251 ; $ws.Cells.Item(9, 8).Value = 16.09209692505834 ; $ws.Cells.Item(9, 10).Value = 9.34309950063523 ; $ws.Cells.Item(9, 11).Value = 10.50357043454236 ; $ws.Cells.Item(9, 15).Value = 25.4540399252061
$ws.Cells.Item(10, 2).Value = 11.85174779217004 ; $ws.Cells.Item(10, 3).Value = 5.481734542846938 ; $ws.Cells.Item(10, 4).Value = 15.26309238434761 ; $ws.Cells.Item(10, 5).Value = 16.6827720408696 ; $ws.Cells.Item(10, 7).Value = 35.84384551858346 ; $ws.Cells.Item(10, 8).Value = 16.01667177628809 ; $ws.Cells.Item(10, 10).Value = 9.342852732029375 ; $ws.Cells.Item(10, 11).Value = 10.86122581032309 ; $ws.Cells.Item(10, 15).Value = 25.34108482996346
$ws.Cells.Item(11, 2).Value = 12.07174676256857 ; $ws.Cells.Item(11, 3).Value = 5.617619542452086 ; $ws.Cells.Item(11, 4).Value = 15.33140619922412 ; $ws.Cells.Item(11, 5).Value = 16.75277625862773 ; $ws.Cells.Item(11, 7).Value = 35.81621895115292 ; $ws.Cells.Item(11, 8).Value = 15.98509936342068 ; $ws.Cells.Item(11, 10).Value = 9.344188832202468 ; $ws.Cells.Item(11, 11).Value = 11.02119507122504 ; $ws.Cells.Item(11, 15).Value = 25.29539917316008
$ws.Cells.Item(12, 2).Value = 12.15419954445891 ; $ws.Cells.Item(12, 3).Value = 5.668118463252393 ; $ws.Cells.Item(12, 4).Value = 15.35768320676681 ; $ws.Cells.Item(12, 5).Value = 16.7797462824064 ; $ws.Cells.Item(12, 7).Value = 35.80740560261717 ; $ws.Cells.Item(12, 8).Value = 15.973537815591 ; $ws.Cells.Item(12, 10).Value = 9.344902290537087 ; $ws.Cells.Item(12, 11).Value = 11.08132386753079 ; $ws.Cells.Item(12, 15).Value = 25.27892114979549
$ws.Cells.Item(13, 2).Value = 12.13648110183564 ; $ws.Cells.Item(13, 3).Value = 5.657285532044016 ; $ws.Cells.Item(13, 4).Value = 15.35200612326302 ; $ws.Cells.Item(13, 5).Value = 16.77391759218662 ; $ws.Cells.Item(13, 7).Value = 35.80923031026455 ; $ws.Cells.Item(13, 8).Value = 15.97601026156499 ; $ws.Cells.Item(13, 10).Value = 9.344739418838964 ; $ws.Cells.Item(13, 11).Value = 11.06839483090719 ; $ws.Cells.Item(13, 15).Value = 25.28243337921093
$ws.Cells.Item(14, 2).Value = 12.07854778049089 ; $ws.Cells.Item(14, 3).Value = 5.621793429499935 ; $ws.Cells.Item(14, 4).Value = 15.3335599500724 ; $ws.Cells.Item(14, 5).Value = 16.75498596394055 ; $ws.Cells.Item(14, 7).Value = 35.81546080269408 ; $ws.Cells.Item(14, 8).Value = 15.98414028450819 ; $ws.Cells.Item(14, 10).Value = 9.344243374842378 ; $ws.Cells.Item(14, 11).Value = 11.02615118938968 ; $ws.Cells.Item(14, 15).Value = 25.29402702389631
$ws.Cells.Item(15, 2).Value = 12.04294827235428 ; $ws.Cells.Item(15, 3).Value = 5.599928159671877 ; $ws.Cells.Item(15, 4).Value = 15.32231374083409 ; $ws.Cells.Item(15, 5).Value = 16.74344928970431 ; $ws.Cells.Item(15, 7).Value = 35.81949199544274 ; $ws.Cells.Item(15, 8).Value = 15.98917151168455 ; $ws.Cells.Item(15, 10).Value = 9.343966532642254 ; $ws.Cells.Item(15, 11).Value = 11.00021579681795 ; $ws.Cells.Item(15, 15).Value = 25.30123561323477
$ws.Cells.Item(16, 2).Value = 11.83725463223602 ; $ws.Cells.Item(16, 3).Value = 5.472721243626664 ; $ws.Cells.Item(16, 4).Value = 15.25868599300845 ; $ws.Cells.Item(16, 5).Value = 16.67826257656902 ; $ws.Cells.Item(16, 7).Value = 35.84588127247433 ; $ws.Cells.Item(16, 8).Value = 16.01879025435071 ; $ws.Cells.Item(16, 10).Value = 9.342794496110999 ; $ws.Cells.Item(16, 11).Value = 10.85071168228959 ; $ws.Cells.Item(16, 15).Value = 25.34418539466526
$ws.Cells.Item(17, 2).Value = 11.70962462585712 ; $ws.Cells.Item(17, 3).Value = 5.393000708563831 ; $ws.Cells.Item(17, 4).Value = 15.2203983084775 ; $ws.Cells.Item(17, 5).Value = 16.63911265674324 ; $ws.Cells.Item(17, 7).Value = 35.86499886441887 ; $ws.Cells.Item(17, 8).Value = 16.03766216268626 ; $ws.Cells.Item(17, 10).Value = 9.342445934114659 ; $ws.Cells.Item(17, 11).Value = 10.75825592164389 ; $ws.Cells.Item(17, 15).Value = 25.37199512623465
$ws.Cells.Item(18, 2).Value = 11.63571011642519 ; $ws.Cells.Item(18, 3).Value = 5.346536356876512 ; $ws.Cells.Item(18, 4).Value = 15.19865597929975 ; $ws.Cells.Item(18, 5).Value = 16.61690895270867 ; $ws.Cells.Item(18, 7).Value = 35.87706876831448 ; $ws.Cells.Item(18, 8).Value = 16.04877455645112 ; $ws.Cells.Item(18, 10).Value = 9.342381828509154 ; $ws.Cells.Item(18, 11).Value = 10.70482395674124 ; $ws.Cells.Item(18, 15).Value = 25.38852672992869
$ws.Cells.Item(19, 2).Value = 11.61059967370062 ; $ws.Cells.Item(19, 3).Value = 5.330700092485298 ; $ws.Cells.Item(19, 4).Value = 15.19134296982712 ; $ws.Cells.Item(19, 5).Value = 16.60944566100974 ; $ws.Cells.Item(19, 7).Value = 35.88133965672453 ; $ws.Cells.Item(19, 8).Value = 16.05258128254066 ; $ws.Cells.Item(19, 10).Value = 9.342383568273505 ; $ws.Cells.Item(19, 11).Value = 10.68669096078064 ; $ws.Cells.Item(19, 15).Value = 25.39421605110195
$ws.Cells.Item(20, 2).Value = 11.72326392175558 ; $ws.Cells.Item(20, 3).Value = 5.401550527010176 ; $ws.Cells.Item(20, 4).Value = 15.22444527177795 ; $ws.Cells.Item(20, 5).Value = 16.64324782865257 ; $ws.Cells.Item(20, 7).Value = 35.86285255841634 ; $ws.Cells.Item(20, 8).Value = 16.03562653460724 ; $ws.Cells.Item(20, 10).Value = 9.342468929762267 ; $ws.Cells.Item(20, 11).Value = 10.76812468179863 ; $ws.Cells.Item(20, 15).Value = 25.36897921975925
$ws.Cells.Item(21, 2).Value = 12.09558802049193 ; $ws.Cells.Item(21, 3).Value = 5.632244474335225 ; $ws.Cells.Item(21, 4).Value = 15.33896710800803 ; $ws.Cells.Item(21, 5).Value = 16.76053427581952 ; $ws.Cells.Item(21, 7).Value = 35.81358597349172 ; $ws.Cells.Item(21, 8).Value = 15.98174159699251 ; $ws.Cells.Item(21, 10).Value = 9.344383449907843 ; $ws.Cells.Item(21, 11).Value = 11.03857175696971 ; $ws.Cells.Item(21, 15).Value = 25.29059935527412
$ws.Cells.Item(22, 2).Value = 12.33390239654133 ; $ws.Cells.Item(22, 3).Value = 5.777425736660908 ; $ws.Cells.Item(22, 4).Value = 15.416183501372 ; $ws.Cells.Item(22, 5).Value = 16.83986581321584 ; $ws.Cells.Item(22, 7).Value = 35.79099702435683 ; $ws.Cells.Item(22, 8).Value = 15.94882278998297 ; $ws.Cells.Item(22, 10).Value = 9.34684378140604 ; $ws.Cells.Item(22, 11).Value = 11.21268948351721 ; $ws.Cells.Item(22, 15).Value = 25.24416686467217
$ws.Cells.Item(23, 2).Value = 12.20719288067337 ; $ws.Cells.Item(23, 3).Value = 5.700457767526566 ; $ws.Cells.Item(23, 4).Value = 15.37476083151813 ; $ws.Cells.Item(23, 5).Value = 16.79728603057172 ; $ws.Cells.Item(23, 7).Value = 35.80217187886115 ; $ws.Cells.Item(23, 8).Value = 15.96618176855056 ; $ws.Cells.Item(23, 10).Value = 9.345420302870162 ; $ws.Cells.Item(23, 11).Value = 11.12001811435154 ; $ws.Cells.Item(23, 15).Value = 25.26850929153606
$ws.Cells.Item(24, 2).Value = 11.71709926682068 ; $ws.Cells.Item(24, 3).Value = 5.397687120781792 ; $ws.Cells.Item(24, 4).Value = 15.22261479803953 ; $ws.Cells.Item(24, 5).Value = 16.64137736886431 ; $ws.Cells.Item(24, 7).Value = 35.86381954300781 ; $ws.Cells.Item(24, 8).Value = 16.03654602383025 ; $ws.Cells.Item(24, 10).Value = 9.342458108851961 ; $ws.Cells.Item(24, 11).Value = 10.76366387671795 ; $ws.Cells.Item(24, 15).Value = 25.37034101853634
$ws.Cells.Item(25, 2).Value = 11.16952175048219 ; $ws.Cells.Item(25, 3).Value = 5.047857456570708 ; $ws.Cells.Item(25, 4).Value = 15.06925904418701 ; $ws.Cells.Item(25, 5).Value = 16.48528313024242 ; $ws.Cells.Item(25, 7).Value = 35.97133841881593 ; $ws.Cells.Item(25, 8).Value = 16.12229889697156 ; $ws.Cells.Item(25, 10).Value = 9.344479039149633 ; $ws.Cells.Item(25, 11).Value = 10.36982746984862 ; $ws.Cells.Item(25, 15).Value = 25.50067029500297
